# Stainless_Steel_Pan_Head_Screws_with_External-Tooth_Lock_Washer.xlsx
#
# Two new rows are inserted above the existing data:
#   - new row 1: a purely numeric header-index row (0..13), replacing the
#     old text header row (keeps the old s="1" header style since we only
#     change the *value*, not the style, of the existing row-1 cells).
#   - new row 2: a blank spacer row except for E2 = "Washer".
#   - the old header row (old row 1) becomes new row 3, but with the
#     K/M/N cells cleared out.
#   - all old data rows (old rows 2..49) shift straight down by two, to
#     new rows 4..51, with identical content.
#
# Because this runtime auto-coerces numeric-looking / comma-having text
# (e.g. "50", "10", "70,000") into real numbers on a plain `.Value =`
# assignment, every *text* write below forces the cell to Text format
# first and then resets the cell style back to "Normal" afterwards so we
# don't leave a stray NumberFormat/quotePrefix style behind (the target
# file has no explicit style on any of these cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 14   # A..N
$lastOldRow = 49  # original sheet was A1:N49

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    if ($text -eq $null) {
        $cell.Value = ""
    } else {
        $cell.NumberFormat = "@"
        $cell.Value = $text
        $cell.Style = "Normal"
    }
}

# 1) Shift all existing rows down by two, bottom-up so we never clobber a
#    row before it has been copied.
for ($oldRow = $lastOldRow; $oldRow -ge 1; $oldRow--) {
    $newRow = $oldRow + 2
    for ($col = 1; $col -le $lastCol; $col++) {
        $text = $ws.Cells.Item($oldRow, $col).Text
        Set-TextCell $newRow $col $text
    }
}

# 2) Old header row is now row 3; clear out K3, M3, N3 per the new layout.
Set-TextCell 3 11 $null
Set-TextCell 3 13 $null
Set-TextCell 3 14 $null

# 3) New row 2 is a spacer row: blank everywhere except E2 = "Washer".
for ($col = 1; $col -le $lastCol; $col++) {
    Set-TextCell 2 $col $null
}
Set-TextCell 2 5 "Washer"

# 4) New row 1 becomes a plain numeric index row 0..13 (keeps the
#    existing bold/bordered header style since we only touch .Value).
for ($col = 1; $col -le $lastCol; $col++) {
    $ws.Cells.Item(1, $col).Value = $col - 1
}
